$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Test-Cases" table originally listed the scenarios in this order:
#   TestScenario_3 (rows 2-5), TestScenario_4 (rows 6-8),
#   TestScenario_1 (rows 9-12), TestScenario_2 (rows 13-14),
#   TestScenario_3 again (rows 15-18)
# This commit reorders the scenario blocks so TestScenario_1 and
# TestScenario_2 come first, followed by the two TestScenario_3 blocks and
# finally TestScenario_4 - i.e. rows 2-8 are moved to the end of the table
# (after row 18). The Cut/Insert COM verbs are unreliable against a sheet
# that hosts a ListObject in this runtime, so the final values are written
# back directly, cell by cell, which yields an identical end result.

$newValues = @(
    @("TestScenario_1", "TestScenario_1.TestCase_1", "New Account", "User Needs to Login to Salesforce, from the browser with correct credentials", "", "Step 1", "Click on the Account tab, and click on New button", "User should be navigated to the New  Account Page", "", ""),
    @("", "", "", "", "Valid value for required field Account Name ", "Step 2", "Input valid value in the  Account Name field.", "User should be able to input value for the Account Name field.", "", ""),
    @("", "", "", "", "Valid value for required field  ", "Step 3", "Input valid value in the   field.", "User should be able to input value for the  field.", "", ""),
    @("", "", "", "", "", "Step 4", "Click on Save button to save Account with fields", "User should be able to validate that a New Account is created", "", ""),
    @("TestScenario_2", "TestScenario_2.TestCase_1", "View Account", "User Needs to Login to Salesforce, from the browser with correct credentials", "", "Step 1", "Click on the Account tab,  and select a Account ", "User should be navigated to the Account Page", "", ""),
    @("", "", "", "", "", "Step 2", "Click on Account name to View the Details", "User should be able to view the Account Details", "", ""),
    @("TestScenario_3", "TestScenario_3.TestCase_1", "Edit Account", "User Needs to Login to Salesforce, from the browser with correct credentials", "", "Step 1", "Click on the Account tab,  and click on existing Account to modify", "User is navigated to the Account Details page", "", ""),
    @("", "", "", "", "Valid value for required field Account Name ", "Step 2", "Input valid value in the  Account Name field.", "User should be able to input value for the Account Name field.", "", ""),
    @("", "", "", "", "Valid value for required field  ", "Step 3", "Input valid value in the   field.", "User should be able to input value for the  field.", "", ""),
    @("", "", "", "", "", "Step 4", "Click on Save button to save Account with fields", "User should be able to validate that the Account is edited", "", ""),
    @("TestScenario_3", "TestScenario_3.TestCase_1", "Edit Account", "User Needs to Login to Salesforce, from the browser with correct credentials", "", "Step 1", "Click on the Account tab,  and click on existing Account to modify", "User is navigated to the Account Details page", "", ""),
    @("", "", "", "", "Valid value for required field Account Name ", "Step 2", "Input valid value in the  Account Name field.", "User should be able to input value for the Account Name field.", "", ""),
    @("", "", "", "", "Valid value for required field  ", "Step 3", "Input valid value in the   field.", "User should be able to input value for the  field.", "", ""),
    @("", "", "", "", "", "Step 4", "Click on Save button to save Account with fields", "User should be able to validate that the Account is edited", "", ""),
    @("TestScenario_4", "TestScenario_4.TestCase_1", "Delete Account", "User Needs to Login to Salesforce, from the browser with correct credentials", "", "Step 1", "Click on the Account tab,  and select the existing  Account to delete", "User is navigated to the Account Details page", "", ""),
    @("", "", "", "", "", "Step 2", "Click on to the Delete to Delete the Account", "User should be able to validate that a pop-up is displayed asking for confirmation to delete the Account", "", ""),
    @("", "", "", "", "", "Step 3", "Click on Confirm / OK to delete the  Account", "User should be able to validate the Account is deleted", "", "")
)

$startRow = 2
for ($i = 0; $i -lt $newValues.Count; $i++) {
    $rowValues = $newValues[$i]
    $targetRow = $startRow + $i
    for ($col = 0; $col -lt $rowValues.Count; $col++) {
        $ws.Cells.Item($targetRow, $col + 1).Value = $rowValues[$col]
    }
}
